$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write order matters for shared-string table layout, mirroring the
# original authoring sequence: B2 before B1, then column-major C, D, E.
$ws.Range("B2").Value = "English"
$ws.Range("B1").Value = "Language Options_Australia"
$ws.Range("C1").Value = "Language Options_Malaysia"
$ws.Range("C2").Value = "English;Bahasa Melayu"
$ws.Range("D1").Value = "Language Options_Philippines"
$ws.Range("D2").Value = "English"
$ws.Range("E1").Value = "Language Options_Thailand"
$ws.Range("E2").Value = "อังกฤษ;ไทย"

# Column widths to match target layout (Excel "best fit" autosize)
$ws.Columns.Item(1).ColumnWidth = 32.92
$ws.Columns.Item(2).ColumnWidth = 23.42
$ws.Columns.Item(3).ColumnWidth = 23.42
$ws.Columns.Item(4).ColumnWidth = 25.09
$ws.Columns.Item(5).ColumnWidth = 23.25

# Selection to match target view
$ws.Range("D11").Select()
